# Árjegyzék (eddigi).xlsx - update pricing figures, remove the
# "Vezetékes telefon" (landline phone) line item, and label the sheet
# with the pricing date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Router quantity dropped from 5 to 3 (formula Q2 recalculates automatically)
$ws.Range("D2").Value = 3

# Bővitő kártya (expansion card) quantity dropped from 6 to 4 (Q4 recalculates)
$ws.Range("D4").Value = 4

# Remove the "Vezetékes telefon" row's data (row 8) entirely, keeping the
# (now empty) row/cell formatting in place, and let the row height revert
# to the sheet default now that the wrapped text is gone.
$ws.Range("A8:E8").ClearContents()
$ws.Rows.Item(8).AutoFit()

# Label the price list with the month the prices were captured.
$ws.Range("A13").Value = "2022. jan. árak"

# Update the saved selection to match the author's last position.
$ws.Range("B14").Select()
